# Update GG Expenditure (U), GG Revenue (V) and GG Balance (W) figures
# for rows 105-123 on the "Quarterly Data Levels" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => U, V, W values
$updates = @{
    105 = @(27499, 25628, -1871)
    106 = @(25558, 19568, -5989)
    107 = @(26016, 23816, -2200)
    108 = @(26576, 24028, -2549)
    109 = @(27381, 32256, 4875)
    110 = @(24686, 24337, -350)
    111 = @(26097, 27634, 1537)
    112 = @(26957, 28777, 1820)
    113 = @(29863, 35091, 5228)
    114 = @(27218, 27625, 407)
    115 = @(28343, 29661, 1318)
    116 = @(29178, 28840, -338)
    117 = @(31882, 37654, 5772)
    118 = @(29693, 28896, -797)
    119 = @(30601, 33532, 2930)
    120 = @(31210, 46615, 15406)
    121 = @(34494, 39560, 5066)
    122 = @(30389, 30909, 521)
    123 = @(32817, 34827, 2010)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("U$row").Value = $vals[0]
    $ws.Range("V$row").Value = $vals[1]
    $ws.Range("W$row").Value = $vals[2]
}
